# Apply weekly fruit/vegetable price update: rows 2-14 get their
# Fecha / Volumen / Precio minimo / Precio maximo / Precio promedio ponderado /
# Origen / Precio $/Kg values redistributed among the existing rows
# (a reshuffle of the weekly records), per the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that get reshuffled between rows.
$cols = @("D", "J", "K", "L", "M", "O", "P")

# Snapshot the "before" values for each affected column, rows 2-14.
$snapshot = @{}
foreach ($col in $cols) {
    $snapshot[$col] = @{}
    for ($r = 2; $r -le 14; $r++) {
        $snapshot[$col][$r] = $ws.Range("$col$r").Value2
    }
}

# Mapping: destination row -> source row (which original row's values now
# live in the destination row).
$mapping = @{
    2  = 8
    3  = 11
    4  = 5
    5  = 6
    6  = 7
    7  = 9
    8  = 13
    9  = 2
    10 = 14
    11 = 4
    12 = 10
    13 = 3
    14 = 12
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $snapshot[$col][$srcRow]
    }
}
